# Update "Pouya Finance" yearly income-statement database for
# saroum/siman (rial.xlsx): roll the 5-period window forward by one
# year (drop the oldest period, shift the remaining four left, and
# append the freshly reported period), refresh the "published on"
# dates, and clear the old placeholder "-" text-value in the
# impairment-expense row (now a proper numeric 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# ---------------------------------------------------------------
# Header row 8: fiscal-period captions (12-month period ended ...)
# ---------------------------------------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# ---------------------------------------------------------------
# Header row 9: publish dates
# ---------------------------------------------------------------
$ws.Range("D9").Value = "1399-01-24 (7)"
$ws.Range("E9").Value = "1400-02-04 (7)"
$ws.Range("F9").Value = "1401-01-31 (8)"
$ws.Range("G9").Value = "1402-01-30 (9)"
$ws.Range("H9").Value = "1402-01-30 (2)"

# ---------------------------------------------------------------
# Row 11: فروش (Sales)
# ---------------------------------------------------------------
$ws.Range("D11").Value = 1700836
$ws.Range("E11").Value = 2251460
$ws.Range("F11").Value = 4346023
$ws.Range("G11").Value = 7506331
$ws.Range("H11").Value = 10215910

# ---------------------------------------------------------------
# Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold)
# ---------------------------------------------------------------
$ws.Range("D12").Value = -1118353
$ws.Range("E12").Value = -1310326
$ws.Range("F12").Value = -2111994
$ws.Range("G12").Value = -3255376
$ws.Range("H12").Value = -4625913

# ---------------------------------------------------------------
# Row 13: سود (زیان) ناخالص (Gross profit)
# ---------------------------------------------------------------
$ws.Range("D13").Value = 582483
$ws.Range("E13").Value = 941134
$ws.Range("F13").Value = 2234029
$ws.Range("G13").Value = 4250955
$ws.Range("H13").Value = 5589997

# ---------------------------------------------------------------
# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
# ---------------------------------------------------------------
$ws.Range("D14").Value = -75595
$ws.Range("E14").Value = -63127
$ws.Range("F14").Value = -85258
$ws.Range("G14").Value = -285043
$ws.Range("H14").Value = -495619

# ---------------------------------------------------------------
# Row 15: هزینه کاهش ارزش دریافتنی‌ها (Impairment expense)
# D15 used to be the text "-" and is now a real numeric zero.
# ---------------------------------------------------------------
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# ---------------------------------------------------------------
# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other op. inc/exp)
# ---------------------------------------------------------------
$ws.Range("D16").Value = 19996
$ws.Range("E16").Value = 41899
$ws.Range("F16").Value = 50256
$ws.Range("G16").Value = 3738
$ws.Range("H16").Value = 368411

# ---------------------------------------------------------------
# Row 17: سود (زیان) عملیاتی (Operating profit)
# ---------------------------------------------------------------
$ws.Range("D17").Value = 526884
$ws.Range("E17").Value = 919906
$ws.Range("F17").Value = 2199027
$ws.Range("G17").Value = 3969650
$ws.Range("H17").Value = 5462789

# ---------------------------------------------------------------
# Row 18: هزینه های مالی (Finance costs)
# ---------------------------------------------------------------
$ws.Range("D18").Value = -2370
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = -5328
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = -30132

# ---------------------------------------------------------------
# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-op. inc/exp)
# ---------------------------------------------------------------
$ws.Range("D19").Value = 20232
$ws.Range("E19").Value = 68671
$ws.Range("F19").Value = 592606
$ws.Range("G19").Value = 969608
$ws.Range("H19").Value = 598855

# ---------------------------------------------------------------
# Row 20: سود خالص عملیات در حال تداوم قبل از مالیات (Profit before tax)
# ---------------------------------------------------------------
$ws.Range("D20").Value = 544746
$ws.Range("E20").Value = 988577
$ws.Range("F20").Value = 2786305
$ws.Range("G20").Value = 4939258
$ws.Range("H20").Value = 6031512

# ---------------------------------------------------------------
# Row 21: مالیات (Tax)
# ---------------------------------------------------------------
$ws.Range("D21").Value = -84040
$ws.Range("E21").Value = -71900
$ws.Range("F21").Value = -203480
$ws.Range("G21").Value = -439928
$ws.Range("H21").Value = -488098

# ---------------------------------------------------------------
# Row 22: سود خالص عملیات در حال تداوم (Net profit, continuing ops)
# ---------------------------------------------------------------
$ws.Range("D22").Value = 460706
$ws.Range("E22").Value = 916677
$ws.Range("F22").Value = 2582825
$ws.Range("G22").Value = 4499330
$ws.Range("H22").Value = 5543414

# ---------------------------------------------------------------
# Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (Discontinued ops)
# ---------------------------------------------------------------
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 6448
$ws.Range("H23").Value = 5543

# ---------------------------------------------------------------
# Row 24: سود (زیان) خالص (Net profit)
# ---------------------------------------------------------------
$ws.Range("D24").Value = 460706
$ws.Range("E24").Value = 916677
$ws.Range("F24").Value = 2582825
$ws.Range("G24").Value = 4499330
$ws.Range("H24").Value = 5543414

# ---------------------------------------------------------------
# Row 25: سود هر سهم پس از کسر مالیات (EPS after tax)
# ---------------------------------------------------------------
$ws.Range("D25").Value = 658
$ws.Range("E25").Value = 1310
$ws.Range("F25").Value = 3690
$ws.Range("G25").Value = 6428
$ws.Range("H25").Value = 5543

# ---------------------------------------------------------------
# Row 26: سرمایه (Capital)
# ---------------------------------------------------------------
$ws.Range("D26").Value = 700000
$ws.Range("E26").Value = 700000
$ws.Range("F26").Value = 700000
$ws.Range("G26").Value = 700000
$ws.Range("H26").Value = 1000000

# ---------------------------------------------------------------
# Row 27: سود هر سهم بر اساس آخرین سرمایه (EPS on latest capital)
# ---------------------------------------------------------------
$ws.Range("D27").Value = 461
$ws.Range("E27").Value = 917
$ws.Range("F27").Value = 2583
$ws.Range("G27").Value = 4499
$ws.Range("H27").Value = 5543
